$wb = $excel.ActiveWorkbook

# Select the Collections worksheet (it had a blank first row that needs removing)
$ws = $wb.Worksheets.Item("Collections")
$ws.Activate()

# Remove the initial blank row (row 1), shifting the data (rows 2-4) up by one
$ws.Rows.Item(1).Delete()

# Reflect the selection left after deleting row 1 in the Collections sheet
$ws.Range("B3").Select()
